$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Username values in the shared strings (adds an extra trailing digit "8")
$ws.Range("E2").Value = "anasule0012345678"
$ws.Range("E3").Value = "blakenailya0012345678"
$ws.Range("E4").Value = "mikeaj0012345678"

# Widen columns D (Photograph) and E (Username)
$ws.Columns.Item(4).ColumnWidth = 39.5
$ws.Columns.Item(5).ColumnWidth = 27

# Move the active selection from F10 to G11
[void]$ws.Range("G11").Select()
